$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data per diff
$ws.Range("D2").Value = "26.606.16"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "1.584.00"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.249"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.53"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0833"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "1.805.91"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").Value = "1.579.21"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.53"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "26.607.65"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "207.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.69%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.30"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.41"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("E28").Value = "  -4.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.26"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.677"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +24.91%  "
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").Value = "1.318.37"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.50"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.34"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.784"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "1.718.56"
$ws.Range("E45").Value = "  -2.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.828"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0990"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0506"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.49"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.93%  "
